# "Read only and input cells demarcation"
#
# - ACE Scoring Calculator sheet: years-of-education (B17) and patient
#   score inputs (B18:B26) get new values, the rest of the sheet
#   (averages/SD lookups + derived Major/Minor-impairment calcs) simply
#   ripples via formula recalculation.
# - Those same input cells are unlocked and the worksheet is protected so
#   everything else becomes read-only, with selection restricted to the
#   unlocked (input) cells only.
# - Window/selection state is nudged on both visible sheets to reflect
#   where the author was last working.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("ACE Scoring Calculator")
$ws2 = $wb.Worksheets.Item("ACE scoring")

# --- Update the patient input values -------------------------------------
$ws1.Range("B17").Value = 5

$ws1.Range("B19").Value = 1
$ws1.Range("B20").Value = 2
$ws1.Range("B21").Value = 8
$ws1.Range("B22").Value = 2
$ws1.Range("B23").Value = 7
$ws1.Range("B24").Value = 9
$ws1.Range("B25").Value = 12
$ws1.Range("B26").Value = 4

# --- Demarcate read-only vs. input cells ----------------------------------
# Unlock the cells that are meant to stay editable (years of education +
# the patient score column), then protect the sheet so every other
# (formula/read-only) cell becomes locked, and restrict selection to the
# unlocked cells only.
$ws1.Range("B17").Locked = $false
$ws1.Range("B18:B26").Locked = $false
$ws1.EnableSelection = 1
$ws1.Protect()

# --- Window / selection state ---------------------------------------------
# Update the (currently inactive) "ACE scoring" sheet's remembered
# selection/scroll position first so activating it isn't the final step.
[void]$ws2.Range("A25").Select()
[void]$ws2.Range("B34:B42").Select()

# Leave "ACE Scoring Calculator" as the active tab with its own
# updated scroll position/selection.
[void]$ws1.Activate()
[void]$ws1.Range("A10").Select()
[void]$ws1.Range("B18").Select()
